# Updated cryptos list on Thu Feb  8 15:57:41 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.369.41"
$ws.Range("E2").Value = "  +5.30%  "
$ws.Range("D3").Value = "2.452.04"
$ws.Range("E3").Value = "  +3.51%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.03"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +5.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.62"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +8.71%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.518"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("E9").Value = "  +10.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.06"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +4.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0807"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +2.24%  "
$ws.Range("E12").Value = "  -2.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.65"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  +1.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.05"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +3.59%  "
$ws.Range("D15").Value = "2.836.01"
$ws.Range("E15").Value = "  +3.69%  "
$ws.Range("D16").Value = "2.451.07"
$ws.Range("E16").Value = "  +3.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.841"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  +4.70%  "
$ws.Range("D18").Value = "45.247.01"
$ws.Range("E18").Value = "  +4.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.41"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +3.49%  "
$ws.Range("E20").Value = "  +1.23%  "
$ws.Range("E21").Value = "  +4.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.25"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +1.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "244.45"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +3.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.30"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +2.79%  "
$ws.Range("E25").Value = "  +3.34%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.54"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +4.37%  "
$ws.Range("E28").Value = "  -7.22%  "
$ws.Range("E29").Value = "  +2.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.96"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +6.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "49.83"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  +3.86%  "
$ws.Range("E32").Value = "  +14.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.45"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +14.40%  "
$ws.Range("E34").Value = "  +4.05%  "
$ws.Range("E35").Value = "  +0.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0768"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  +3.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.91"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  +4.93%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.52"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +4.33%  "
$ws.Range("E39").Value = "  +0.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "125.40"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -2.18%  "
$ws.Range("E41").Value = "  +2.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.20"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -2.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.31"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +0.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0292"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +5.05%  "
$ws.Range("D45").Value = "1.947.89"
$ws.Range("E45").Value = "  +0.94%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.99"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +8.05%  "
$ws.Range("E47").Value = "  -1.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.26"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +0.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.79"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +16.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "76.47"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +6.67%  "
$ws.Range("E51").Value = "  +4.12%  "
